$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, row 1 (next to existing "sum" header in G1)
$ws.Range("H1").Value = "Save"

# Copy the header formatting (bold, centered, bordered) from the existing G1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Add the numeric "Save" value for the data row
$ws.Range("H2").Value = 1
